$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# The AR2 winter forecast table gained a new first observation (2007,
# date serial 39400) and every y_0 / y_1 forecast value (columns C & E)
# was recomputed. That pushes the existing 17 data rows (rows 2-18) one
# row down (rows 3-19) while keeping their date / y_0 / y_1 "actuals"
# (columns A, B, D) intact.
# ---------------------------------------------------------------------

# 1) Shift the existing actuals (A, B, D) down by one row, bottom-up so
#    we never overwrite a row before reading it. Value2 is used to read
#    back raw numeric values (avoids any date/variant coercion on write).
for ($r = 18; $r -ge 2; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst").Value = $ws.Range("A$r").Value2()
    $ws.Range("B$dst").Value = $ws.Range("B$r").Value2()
    $ws.Range("D$dst").Value = $ws.Range("D$r").Value2()
}

# 2) Give the brand-new row 19's date cell (A19) the same number format
#    as the rest of column A, without introducing a new style entry.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# 3) Fill in the new leading observation (row 2): 2007-01-01-ish serial
#    date 39400, year 2007, forecasted y_0 value, and the following
#    year label for y_1 (2008). This row never had a y_1 forecast (E).
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("D2").Value = 2008

# 4) Re-populate the recomputed forecast columns C (y_0_forecast) and E
#    (y_1_forecast) for every row of the table.
$ws.Range("C2").Value = 1.75539628881467
$ws.Range("C3").Value = 2.213911448916162
$ws.Range("C4").Value = 2.533533936850563
$ws.Range("C5").Value = 2.088987486264915
$ws.Range("E5").Value = 1.485473821631844
$ws.Range("C6").Value = 1.212544822741002
$ws.Range("E6").Value = 1.799394172339341
$ws.Range("C7").Value = 1.196776590518644
$ws.Range("E7").Value = 1.2151583353186
$ws.Range("C8").Value = 0.4712609263772594
$ws.Range("E8").Value = 1.107727073902187
$ws.Range("C9").Value = 0.8783377572271434
$ws.Range("E9").Value = 1.612081704302182
$ws.Range("C10").Value = 2.29066283401107
$ws.Range("E10").Value = 2.221748592150097
$ws.Range("C11").Value = 4.109890522944348
$ws.Range("E11").Value = 2.932944072183674
$ws.Range("C12").Value = 1.336316831462692
$ws.Range("E12").Value = 1.104283769064729
$ws.Range("C13").Value = 1.197912858979611
$ws.Range("E13").Value = 1.649865498505276
$ws.Range("C14").Value = 1.727537197898665
$ws.Range("E14").Value = 2.284828905445169
$ws.Range("C15").Value = 3.647228437274408
$ws.Range("E15").Value = 3.474365686630398
$ws.Range("C16").Value = 2.777797690741424
$ws.Range("E16").Value = 1.742844348069261
$ws.Range("C17").Value = 0.6994919452575576
$ws.Range("E17").Value = 0.5651273241891186
$ws.Range("C18").Value = -1.432689847121871
$ws.Range("E18").Value = 0.4518870186319468
$ws.Range("C19").Value = 2.033479419175133
$ws.Range("E19").Value = 1.959987726090251
